$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the bold paragraph "Розглядається симетричний, асиметричний
#    та змішаний варіанти." and the empty numbered (numId=6) paragraph
#    that used to follow the "В загальному випадку..." paragraph.
#    We collect the target paragraph indices first (content-matched,
#    not position-matched) and then delete from the bottom up so that
#    earlier indices stay valid while we work.
# ------------------------------------------------------------------

$count = $d.Paragraphs.Count
$deleteIdx = @()

for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -like "*Розглядається симетричний*") {
        $deleteIdx += $i
    }
    if ($t -like "*В загальному випадку, асиметрична задача комівояжера*") {
        $deleteIdx += ($i + 1)
    }
}

$deleteIdx = $deleteIdx | Sort-Object -Descending
foreach ($idx in $deleteIdx) {
    $d.Paragraphs.Item($idx).Range.Delete()
}

# ------------------------------------------------------------------
# 2) Flip the language of the two empty paragraphs that currently
#    carry en-US on their paragraph mark over to ru-RU (identified by
#    being empty AND currently en-US, so the title run "NP" - which is
#    not an empty paragraph - is left untouched).
# ------------------------------------------------------------------

$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.Length -le 1 -and $p.Range.LanguageID -eq "en-US") {
        $p.Range.LanguageID = "ru-RU"
    }
}
